$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B13").Value = "ひゃくえむ。"
$ws.Range("B14").Value = "南海トラフ巨大地震"
$ws.Range("B15").Value = "黒猫と魔女の教室"
$ws.Range("B16").Value = "魔女と傭兵"
$ws.Range("B22").Value = "魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～"
$ws.Range("B23").Value = "味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す"
$ws.Range("B24").Value = "ドラハチ"
$ws.Range("B25").Value = "アルキメデスの大戦"
$ws.Range("B26").Value = "四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～"
$ws.Range("B33").Value = "せいぶつ部の田辺くん"
$ws.Range("B34").Value = "おやすみ ふみさん"
$ws.Range("B35").Value = "阿武ノーマル"
$ws.Range("B38").Value = "イレギュラーズ"
$ws.Range("B39").Value = "Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜"
$ws.Range("B40").Value = "FAIRY TAIL 100 YEARS QUEST"
$ws.Range("B41").Value = "念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！"
$ws.Range("B43").Value = "アオバノバスケ"
$ws.Range("B45").Value = "食糧人類-Starving Anonymous-"
$ws.Range("B51").Value = "なれの果ての僕ら"
$ws.Range("B52").Value = "イジらないで、長瀞さん"
$ws.Range("B53").Value = "Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。"
$ws.Range("B56").Value = "最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～"
$ws.Range("B57").Value = "辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜"
$ws.Range("B58").Value = "いじめるヤバイ奴"
$ws.Range("B59").Value = "ハナバス　苔石花江のバスケ論"
$ws.Range("B60").Value = "デスティニーラバーズ"
$ws.Range("B62").Value = "降り積もれ孤独な死よ"
$ws.Range("B63").Value = "不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～"
$ws.Range("B64").Value = "ブルーロック"
$ws.Range("B65").Value = "追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜"
$ws.Range("B66").Value = "ジュミドロ"
$ws.Range("B67").Value = "リスナーに騙されてダンジョンの最下層から脱出RTAすることになった"
$ws.Range("B68").Value = "剣帝学院の魔眼賢者"
$ws.Range("B69").Value = "人間消失"
$ws.Range("B70").Value = "可愛いだけじゃない式守さん"
$ws.Range("B71").Value = "幼馴染とはラブコメにならない"
$ws.Range("B72").Value = "ヒロインは絶望しました。"
$ws.Range("B73").Value = "魁の花巫女"
$ws.Range("B74").Value = "復讐の教科書"
$ws.Range("B75").Value = "MYS"
$ws.Range("B76").Value = "ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～"
$ws.Range("B77").Value = "異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～"
$ws.Range("B79").Value = "シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～"
$ws.Range("B80").Value = "東京卍リベンジャーズ～場地圭介からの手紙～"
$ws.Range("B81").Value = "我間乱 ―修羅―"
$ws.Range("B82").Value = "インフェクション"
$ws.Range("B83").Value = "DAYS外伝"
$ws.Range("B85").Value = "はっちぽっちぱんち"
$ws.Range("B86").Value = "勇者と呼ばれた後に　―そして無双男は家族を創る―"
$ws.Range("B87").Value = "卒業アルバムの彼女たち"
$ws.Range("B88").Value = "ぼくのアデリア"
$ws.Range("B89").Value = "追放されなかった男　～二度目の人生は土下座から始まりました～"
$ws.Range("B90").Value = "GALAXIAS"
$ws.Range("B91").Value = "ともだちづくり"
$ws.Range("B92").Value = "劣等人の魔剣使い　スキルボードを駆使して最強に至る"
$ws.Range("B93").Value = "お願い、脱がシて。"
$ws.Range("B94").Value = "田んぼで拾った女騎士、田舎で俺の嫁だと思われている"
$ws.Range("B95").Value = "はじめの一歩"
$ws.Range("B96").Value = "「俺、パーティー抜けるわ」が口癖のスキル【縮小】のDランク冒険者、聖女と結婚して勇者パーティーに加入するハメになる"
$ws.Range("B97").Value = "冰剣の魔術師が世界を統べる〜世界最強の魔術師である少年は、魔術学院に入学する〜"
$ws.Range("B98").Value = "彼女、お借りします"
$ws.Range("B99").Value = "それがメイドのカンナです"
$ws.Range("B100").Value = "日本語が話せないロシア人美少女転入生が頼れるのは、多言語マスターの俺1人"
$ws.Range("B101").Value = "ザ・ファブル"
